# 2021 PLANETE SOLIDAIRE - CONTACTS PARTENAIRES
# Reconcile the partner contact list: drop duplicate/obsolete contacts,
# fix name/typo issues, and rebuild the "Nom"/"Prenom" ordering + mail hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove every existing hyperlink first. Row deletion below does not shift
#    the worksheet's cached hyperlink references in this host, so the safest
#    path is to drop them all and re-add the ones we still need afterwards.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2) Delete the rows for contacts that are no longer tracked:
#      row 16 - Ristorante Duomo (placeholder, no contact)
#      row 15 - Signe Toque (placeholder, no contact)
#      row 14 - Esprit de voyages (placeholder, no contact)
#      row 8  - Mama Shelter / Vincent Brun
#      row 6  - Groupe Pic / Charlotte Terrier
#      row 5  - Groupe Pic / Pauline Monot
#    Delete bottom-up so earlier row numbers stay valid.
# ---------------------------------------------------------------------------
$ws.Rows(16).Delete()
$ws.Rows(15).Delete()
$ws.Rows(14).Delete()
$ws.Rows(8).Delete()
$ws.Rows(6).Delete()
$ws.Rows(5).Delete()

# ---------------------------------------------------------------------------
# 3) Fix remaining rows: correct partner-name typos/casing and swap the
#    "Nom"/"Prenom" values back into their proper columns.
# ---------------------------------------------------------------------------

# Row 3 - Groupe Pic -> Groupe PIC ; David Sinapian is Proprietaire
$ws.Range("A3").Value = "Groupe PIC"
$ws.Range("D3").Value = "Propriétaire"

# Row 4 - Daily Pic -> Daily PIC
$ws.Range("A4").Value = "Daily PIC"

# Row 5 - Mama Shelter / Serge Trigano (Nom/Prenom were swapped)
$ws.Range("B5").Value = "Trigano"
$ws.Range("C5").Value = "Serge"

# Row 6 - Chateau La Brande -> Château LaBrande ; Maude Soulies (Nom/Prenom swapped)
$ws.Range("A6").Value = "Château LaBrande"
$ws.Range("B6").Value = "Souliès"
$ws.Range("C6").Value = "Maude"

# Row 7 - Oui Chef -> Oui Chef! ; Rutger Eysvogel (Nom/Prenom swapped)
$ws.Range("A7").Value = "Oui Chef!"
$ws.Range("B7").Value = "Eysvogel"
$ws.Range("C7").Value = "Rutger"

# Row 8 - Tomeet / Tom Chauvet (Nom/Prenom swapped)
$ws.Range("B8").Value = "Chauvet"
$ws.Range("C8").Value = "Tom"

# Row 9 - Les arcs -> Les Arcs ; Marine Ruas (Nom/Prenom swapped)
$ws.Range("A9").Value = "Les Arcs"
$ws.Range("B9").Value = "Ruas"
$ws.Range("C9").Value = "Marine"

# Row 10 - Village Pro BTP -> Villages Pro BTP ; Stephane -> Stephanie Besnier
$ws.Range("A10").Value = "Villages Pro BTP"
$ws.Range("B10").Value = "Besnier"
$ws.Range("C10").Value = "Stéphanie"

# Row 11 - A la decouverte des langues / Eric Humblot (Nom/Prenom swapped)
$ws.Range("B11").Value = "Humblot"
$ws.Range("C11").Value = "Eric"

# Row 12 - Artiste -> Artiste Dumas ; Francois Dumas (Nom/Prenom swapped)
$ws.Range("A12").Value = "Artiste Dumas"
$ws.Range("B12").Value = "Dumas"
$ws.Range("C12").Value = "François"

# ---------------------------------------------------------------------------
# 4) Re-create the mail hyperlinks on the surviving rows, in the same order
#    they originally appeared, so relationship ids line up the same way.
#    Passing the "mailto:..." text as TextToDisplay caches that string in
#    the hyperlink's display attribute; we then restore the cell's own text
#    (the bare e-mail address) right after.
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:srlauredumas@yahoo.fr", "", "", "mailto:srlauredumas@yahoo.fr")
$ws.Range("E2").Value = "srlauredumas@yahoo.fr"
$ws.Range("E2").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:david.sinapian@groupe-pic.com", "", "", "mailto:david.sinapian@groupe-pic.com")
$ws.Range("E3").Value = "david.sinapian@groupe-pic.com "
$ws.Range("E3").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:lilian.etienne@dailypic.com", "", "", "mailto:lilian.etienne@dailypic.com")
$ws.Range("E4").Value = "lilian.etienne@dailypic.com "
$ws.Range("E4").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:contact@chateau-labrande.fr")
$ws.Range("E6").Value = "contact@chateau-labrande.fr"
$ws.Range("E6").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:rutgereysvogel@gmail.com")
$ws.Range("E7").Value = "rutgereysvogel@gmail.com"
$ws.Range("E7").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:m.ruas@bourgsaintmaurice.fr", "", "", "mailto:m.ruas@bourgsaintmaurice.fr")
$ws.Range("E9").Value = "  m.ruas@bourgsaintmaurice.fr"
$ws.Range("E9").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:tom.chauvet@isg.fr")
$ws.Range("E8").Value = "tom.chauvet@isg.fr"
$ws.Range("E8").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:humblot27@wanadoo.fr")
$ws.Range("E11").Value = "humblot27@wanadoo.fr"
$ws.Range("E11").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:serge@mamashelter.com", "", "", "mailto:serge@mamashelter.com")
$ws.Range("E5").Value = "serge@mamashelter.com"
$ws.Range("E5").Style = "Lien hypertexte"
